$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - D1 becomes a TRUE() formula with the "TRUE"/"FALSE" number format (reuses style 1)
$ws.Range("D1").Formula = "=TRUE()"
$ws.Range("D1").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws.Range("E1").Clear()

# Row 5 - new row mirroring row 1 with a FALSE() formula
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "name2"
$ws.Range("C5").Value = 21.3
$ws.Range("D5").Formula = "=FALSE()"
$ws.Range("D5").NumberFormat = """TRUE"";""TRUE"";""FALSE"""

# Rows 7-9 - list<String> test data
$ws.Range("A7").Value = "n1"
$ws.Range("B7").Value = "n2"
$ws.Range("C7").Value = "n3"
$ws.Range("D7").Value = "n4"

$ws.Range("A8").Value = "n5"
$ws.Range("B8").Value = "n6"
$ws.Range("C8").Value = "n7"
$ws.Range("D8").Value = "n8"

$ws.Range("A9").Value = "n9"
$ws.Range("B9").Value = "n10"
$ws.Range("C9").Value = "n11"
$ws.Range("D9").Value = "n12"

$ws.Range("B7:D9").Font.Name = "Arial"

$ws.Range("C8").Select()
